$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.109.00"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.898.94"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3896"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07866"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9906"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "1.891.67"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.791"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.055"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009929"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "29.113.60"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "2.099.93"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.897"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.877"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09328"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.323"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.160"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05799"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02078"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.701"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5678"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1794"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.744"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.239"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5349"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.849"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.550"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
